$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'22.487.75"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +0.49%  '

$ws.Cells.Item(3, 4).Value = "'1.572.04"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +0.26%  '

$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.01%  '

$ws.Cells.Item(5, 5).Value = '  -0.04%  '

$ws.Cells.Item(6, 4).Value = "'291.27"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.14%  '

$ws.Cells.Item(7, 4).Value = "'0.3711"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -1.49%  '

$ws.Cells.Item(8, 4).Value = "'49.92"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +1.85%  '

$ws.Cells.Item(9, 4).Value = "'0.3378"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.52%  '

$ws.Cells.Item(10, 4).Value = "'1.144"
$ws.Cells.Item(10, 4).Style = "Normal"

$ws.Cells.Item(11, 4).Value = "'0.07536"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.64%  '

$ws.Cells.Item(12, 5).Value = '  -0.06%  '

$ws.Cells.Item(13, 4).Value = "'21.15"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.65%  '

$ws.Cells.Item(14, 4).Value = "'6.015"
$ws.Cells.Item(14, 4).Style = "Normal"

$ws.Cells.Item(15, 4).Value = "'6.953"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.65%  '

$ws.Cells.Item(16, 4).Value = "'1.570.56"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.26%  '

$ws.Cells.Item(17, 4).Value = "'0.00001120"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.60%  '

$ws.Cells.Item(18, 4).Value = "'90.54"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.83%  '

$ws.Cells.Item(19, 4).Value = "'0.06775"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.46%  '

$ws.Cells.Item(20, 5).Value = '  -0.06%  '

$ws.Cells.Item(21, 4).Value = "'6.338"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.10%  '

$ws.Cells.Item(22, 4).Value = "'16.42"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -1.03%  '

$ws.Cells.Item(23, 4).Value = "'12.22"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +2.31%  '

$ws.Cells.Item(24, 4).Value = "'22.480.74"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.48%  '

$ws.Cells.Item(25, 4).Value = "'2.376"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.10%  '

$ws.Cells.Item(26, 4).Value = "'2.615"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.49%  '

$ws.Cells.Item(27, 4).Value = "'20.04"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.74%  '

$ws.Cells.Item(28, 4).Value = "'149.03"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +1.05%  '

$ws.Cells.Item(29, 4).Value = "'5.057"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.82%  '

$ws.Cells.Item(30, 4).Value = "'125.25"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.38%  '

$ws.Cells.Item(31, 4).Value = "'1.747.34"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.42%  '

$ws.Cells.Item(32, 4).Value = "'1.069"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +7.63%  '

$ws.Cells.Item(33, 4).Value = "'6.206"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +2.44%  '

$ws.Cells.Item(34, 5).Value = '  -0.33%  '

$ws.Cells.Item(35, 4).Value = "'9.754"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -3.32%  '

$ws.Cells.Item(36, 4).Value = "'0.08344"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -2.02%  '

$ws.Cells.Item(37, 4).Value = "'1.369"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -3.32%  '

$ws.Cells.Item(38, 4).Value = "'0.02475"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.79%  '

$ws.Cells.Item(39, 4).Value = "'0.2301"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +0.39%  '

$ws.Cells.Item(40, 4).Value = "'0.06537"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +1.51%  '

$ws.Cells.Item(41, 4).Value = "'5.431"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.27%  '

$ws.Cells.Item(42, 4).Value = "'11.28"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.11%  '

$ws.Cells.Item(43, 4).Value = "'0.6208"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -1.81%  '

$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).Value = "'14.10"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.50%  '

$ws.Cells.Item(45, 2).Value = 'Frax'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(45, 4).Value = "'1.001"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.00%  '

$ws.Cells.Item(46, 4).Value = "'3.808"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.18%  '

$ws.Cells.Item(47, 4).Value = "'0.5843"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -1.33%  '

$ws.Cells.Item(48, 4).Value = "'129.22"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +3.74%  '

$ws.Cells.Item(49, 4).Value = "'2.069"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.15%  '

$ws.Cells.Item(50, 5).Value = '  -3.08%  '

$ws.Cells.Item(51, 4).Value = "'0.07318"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.06%  '
